# Navigator Performance Testing.xlsx - FloatCanvas perf numbers updated
# (commit: "fixed some stuff in FloatCanvas and added some optimizations")
#
# The re-measured "FloatCanvas" row's raw elapsed time (B3) dropped from
# 194.18 to 75.9709; its dependent "per-100" formula in C3 recalculates
# automatically. After typing the new value and pressing Enter, the
# selection cursor lands one row below the edited block, on B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the FloatCanvas elapsed-time measurement; C3's shared formula
# (=B3/100) recalculates on its own.
$ws.Range("B3").Value = 75.9709

# Move/collapse the selection to B5, matching the post-edit cursor position.
$ws.Range("B5").Select()
